# Append two new log rows (26 and 27) to the "check_availability" sheet,
# matching the target diff: dimension grows from A1:F25 to A1:F27.
#
# Columns: A=Timestamp, B=Command, C=URL, D=Result, E=Entered Date, F=Entered Time
#
# Columns E and F hold plain date/time-looking text ("2024-09-12", "18:30:00")
# that Excel would otherwise auto-coerce into numeric date/time serials. The
# source file stores every cell (including these) as plain inline text, so we
# force text entry with a leading apostrophe and then reset the cell style
# back to "Normal" so no date/time number format sticks to the new cells
# (keeping them stock, unstyled cells just like the rest of the data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$value)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 26
$ws.Range("A26").Value = "2024-09-12 18:32:12"
$ws.Range("B26").Value = "check_availability"
$ws.Range("C26").Value = "https://www.opentable.com/r/bar-spero-washington/"
$ws.Range("D26").Value = "No availability for the selected date."
Set-TextCell $ws.Range("E26") "2024-09-12"
Set-TextCell $ws.Range("F26") "18:30:00"

# Row 27
$ws.Range("A27").Value = "2024-09-12 18:34:16"
$ws.Range("B27").Value = "check_availability"
$ws.Range("C27").Value = "MOCKURL_https://www.opentable.com/r/bar-spero-washington/"
$ws.Range("D27").Value = "MOCK_No availability for the selected date."
Set-TextCell $ws.Range("E27") "2024-09-12"
Set-TextCell $ws.Range("F27") "18:34:16"
